$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Move the footer (signature) block down from rows 25:26 to rows 32:33
#    so there is room for the expanded data table. Copy formatting+values
#    first while the source rows still hold the original footer content.
#    (Only B:C and H:J are populated in the footer rows; D:G are untouched
#    there, so copy those two sub-ranges rather than the whole B:J block.)
# ---------------------------------------------------------------------------
$ws.Range("B25:C26").Copy($ws.Range("B32:C33"))
$ws.Range("H25:J26").Copy($ws.Range("H32:J33"))

# ---------------------------------------------------------------------------
# 2) Expand the data table from 5 rows (16:20) to 12 rows (16:27).
#    Row 16 carries the "regular" row style, row 20 carries the "last row"
#    (bottom-border) style. Grab the bottom-border style for the new last
#    row (27) before anything else is touched, then clear the stale
#    leftovers (old row 20 values + old footer remnants in rows 25:26)
#    and stamp the "regular" style down across rows 20:26.
# ---------------------------------------------------------------------------
$ws.Range("B20:J20").Copy($ws.Range("B27:J27"))
$ws.Range("B20:J26").ClearContents()
$ws.Range("B16:J16").Copy($ws.Range("B20:J20"))
$ws.Range("B16:J16").Copy($ws.Range("B21:J21"))
$ws.Range("B16:J16").Copy($ws.Range("B22:J22"))
$ws.Range("B16:J16").Copy($ws.Range("B23:J23"))
$ws.Range("B16:J16").Copy($ws.Range("B24:J24"))
$ws.Range("B16:J16").Copy($ws.Range("B25:J25"))
$ws.Range("B16:J16").Copy($ws.Range("B26:J26"))

# ---------------------------------------------------------------------------
# 3) Write the new table contents (new employee + renumbered/extra periods).
# ---------------------------------------------------------------------------
$data = @(
    @(16, "CC", "45451678",   "ANA MARIA GONZALEZ GONZALEZ",   "2507",  89060, 2226500),
    @(17, "CC", "45451678",   "ANA MARIA GONZALEZ GONZALEZ",   "2506",  89060, 2226500),
    @(18, "CC", "45451678",   "ANA MARIA GONZALEZ GONZALEZ",   "2505",  89060, 2226500),
    @(19, "CC", "45451678",   "ANA MARIA GONZALEZ GONZALEZ",   "2504",  89060, 2226500),
    @(20, "CC", "45451678",   "ANA MARIA GONZALEZ GONZALEZ",   "2503",  89060, 2226500),
    @(21, "CC", "45451678",   "ANA MARIA GONZALEZ GONZALEZ",   "2502",  89060, 2226500),
    @(22, "CC", "91541024",   "PEDRO HERLEY RAMIREZ CARDENAS", "1910", 110600, 4369800),
    @(23, "CC", "91541024",   "PEDRO HERLEY RAMIREZ CARDENAS", "1909",  47927, 4369800),
    @(24, "CC", "1050971007", "VICTOR MANUEL MEJIA BABILONIA", "2507",  56940, 1423500),
    @(25, "CC", "1050971007", "VICTOR MANUEL MEJIA BABILONIA", "2506",  56940, 1423500),
    @(26, "CC", "1050971007", "VICTOR MANUEL MEJIA BABILONIA", "2505",  56940, 1423500),
    @(27, "CC", "1050971007", "VICTOR MANUEL MEJIA BABILONIA", "2504",  56940, 1423500)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
}

# ---------------------------------------------------------------------------
# 4) Update the summary/header area: total mora, worker count, period count.
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = "ESTADO DE CUENTA"
$ws.Range("B7").Value = "RAZON SOCIAL:"
$ws.Range("B11").Value = "VALOR MORA"
$ws.Range("E11").Value = 920647
$ws.Range("B13").Value = "Cant. Trabajadores"
$ws.Range("C13").Value = 3
$ws.Range("E13").Value = "Cant. Periodos"
$ws.Range("F13").Value = 8
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"
$ws.Range("J15").Value = "Observaciones"

# ---------------------------------------------------------------------------
# 5) Footer text (kept identical in wording, only its row position changed).
# ---------------------------------------------------------------------------
$ws.Range("B32").Value = "___________________________________"
$ws.Range("H32").Value = "___________________________________"
$ws.Range("B33").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H33").Value = "FIRMA DEL REPRESENTANTE LEGAL"
